$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Labyrinth Monster" raid as row 7, following the existing
# name / story / raid_type / raid_boss_id / raid_monster_ids /
# raid_boss_location_id / corrupted_location_ids /
# scheduled_event_description / item_specialty_reward_type /
# artifact_item_id column layout.
$ws.Range("A7").Value = 'The Labyrinth Monster'
$ws.Range("B7").Value = 'She dances, dressed in white lace and silk. He holds her close, handsome, a prince at her side. They dance in the crowded ball room. She appears, chanting and with magical rage. The Witch plaxces the curse, he screams and roars and the women in white begins to scream. The little girl watches from the shadows as the world burns around her.'
$ws.Range("C7").Value = 'enraged-little-girl'
$ws.Range("D7").Value = 'Enraged and Corrupted Little Girl'
$ws.Range("E7").Value = 'Wolves of the labyrinth,Haunted Maiden in the Maze,Lost child of heartbreak,Labyrinth Monk,Fettered Maiden of Deception,Thieves of Galitonore,Cultist of tomorrow,Scholar of corrupted Alchemy,Cursed Minotaur,Witch of the abyss,Bard of lost songs,Corrupted and vengeful memory,Reaper of the maze,Chains of the walls'
$ws.Range("F7").Value = 'Minotar Port'
$ws.Range("G7").Value = 'Travelers Port,Minotar Port,Town Of Sorrow,Golden Halls'
$ws.Range("H7").Value = 'The Labyrinth Monster Raid brings about the story of The Little Girl who''s mother and father were cursed by an evil witch. Although not everything is as it would seem. Whats mysteries will you unlock as you progress through the story line to take on the Raid Boss?'
$ws.Range("I7").Value = 'Labyrinth Cloth'
$ws.Range("J7").Value = 'Blacksmiths Anvil'

# The new raid_type ("enraged-little-girl") and raid_boss_id ("Enraged and
# Corrupted Little Girl") values are longer than anything previously in
# columns C/D, so Excel's best-fit column sizing widens those two columns.
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 39.166666666666664
